$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency Price (column D) and Volume(1h) (column E) figures
# pulled on Mon Apr  3 04:50:18 UTC 2023. Values that look numeric are
# entered with a leading text-quote (') so they stay plain text (matching
# the original "." thousands-style / zero-padded formatting) instead of
# being auto-coerced into numbers by Excel.

$ws.Range('D2').Value = '27.726.86'
$ws.Range('E2').Value = '  -2.71%  '
$ws.Range('D3').Value = '1.779.86'
$ws.Range('E3').Value = '  -2.18%  '
$ws.Range('D4').Value = "'1.003"
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'310.50"
$ws.Range('E5').Value = '  -1.91%  '
$ws.Range('D6').Value = "'1.003"
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = "'0.5119"
$ws.Range('E7').Value = '  -0.76%  '
$ws.Range('D8').Value = "'0.3787"
$ws.Range('E8').Value = '  -2.37%  '
$ws.Range('D9').Value = "'0.07786"
$ws.Range('E9').Value = '  -7.93%  '
$ws.Range('D10').Value = "'41.19"
$ws.Range('E10').Value = '  -1.55%  '
$ws.Range('D11').Value = "'1.083"
$ws.Range('E11').Value = '  -2.48%  '
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').Value = "'6.202"
$ws.Range('E13').Value = '  -3.50%  '
$ws.Range('D14').Value = "'20.08"
$ws.Range('E14').Value = '  -4.38%  '
$ws.Range('D15').Value = '1.775.13'
$ws.Range('E15').Value = '  -2.54%  '
$ws.Range('D16').Value = "'7.162"
$ws.Range('E16').Value = '  -4.50%  '
$ws.Range('D17').Value = "'91.23"
$ws.Range('E17').Value = '  -1.69%  '
$ws.Range('D18').Value = "'0.00001072"
$ws.Range('E18').Value = '  -5.53%  '
$ws.Range('D19').Value = "'0.06563"
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('E21').Value = '  -4.18%  '
$ws.Range('D22').Value = "'5.912"
$ws.Range('E22').Value = '  -2.76%  '
$ws.Range('D23').Value = '27.779.21'
$ws.Range('E23').Value = '  -2.67%  '
$ws.Range('E24').Value = '  -3.61%  '
$ws.Range('D25').Value = "'2.234"
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('D26').Value = "'159.31"
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  -3.84%  '
$ws.Range('D28').Value = '1.984.22'
$ws.Range('E28').Value = '  -2.31%  '
$ws.Range('D29').Value = "'2.347"
$ws.Range('E29').Value = '  -2.71%  '
$ws.Range('D30').Value = "'125.05"
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('D31').Value = "'0.1070"
$ws.Range('E31').Value = '  -1.45%  '
$ws.Range('D32').Value = "'1.028"
$ws.Range('E32').Value = '  -5.98%  '
$ws.Range('D33').Value = "'3.630"
$ws.Range('E33').Value = '  -1.32%  '
$ws.Range('D34').Value = "'5.470"
$ws.Range('E34').Value = '  -4.61%  '
$ws.Range('D35').Value = "'0.07064"
$ws.Range('E35').Value = '  -5.53%  '
$ws.Range('D36').Value = "'0.02310"
$ws.Range('E36').Value = '  -2.34%  '
$ws.Range('D37').Value = "'8.716"
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('D38').Value = "'0.2119"
$ws.Range('E38').Value = '  -5.11%  '
$ws.Range('D39').Value = "'11.48"
$ws.Range('E39').Value = '  +2.10%  '
$ws.Range('D40').Value = "'5.002"
$ws.Range('E40').Value = '  -3.77%  '
$ws.Range('D41').Value = "'0.6070"
$ws.Range('E41').Value = '  -3.97%  '
$ws.Range('D42').Value = "'1.003"
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('E43').Value = '  -3.95%  '
$ws.Range('D44').Value = "'1.321"
$ws.Range('E44').Value = '  -5.76%  '
$ws.Range('D45').Value = "'13.06"
$ws.Range('E45').Value = '  -3.77%  '
$ws.Range('D46').Value = "'0.5915"
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('D47').Value = "'3.708"
$ws.Range('E47').Value = '  -1.86%  '
$ws.Range('D48').Value = "'127.55"
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').Value = "'1.197"
$ws.Range('E49').Value = '  -0.23%  '
$ws.Range('D50').Value = "'1.891"
$ws.Range('E50').Value = '  -5.06%  '
$ws.Range('D51').Value = "'0.06813"
$ws.Range('E51').Value = '  -2.35%  '
